# Update "想去人数" (F) and "最低票价" (G) figures to the latest scraped
# values, as published to gh-pages.
#
# Sheet "展览" (Worksheets(1)) and "全部类型" (Worksheets(4)) both contain
# the same list of exhibitions, so both get the same updated numbers.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Cell => new value, for the "展览" sheet.
$exhibitionUpdates = @{
    "F3"  = 8578
    "F4"  = 2705
    "F5"  = 954
    "F8"  = 665
    "F9"  = 125
    "F12" = 896
    "F13" = 3703
    "F14" = 266
    "F15" = 149
    "G15" = 19.9
    "F17" = 775
    "F19" = 495
    "F22" = 1346
    "F23" = 1356
    "F24" = 454
    "F27" = 162
    "F28" = 350
    "F29" = 59
    "F30" = 1017
    "F32" = 511
    "F33" = 666
    "F34" = 45
    "F35" = 45
    "F36" = 80
    "F39" = 166
    "F40" = 6
}

# Cell => new value, for the "全部类型" sheet (same data, shifted by one
# row because that sheet has an extra leading entry).
$allTypesUpdates = @{
    "F4"  = 8578
    "F5"  = 2705
    "F6"  = 954
    "F9"  = 665
    "F10" = 125
    "F13" = 896
    "F15" = 3703
    "F16" = 266
    "F17" = 149
    "G17" = 19.9
    "F20" = 775
    "F23" = 495
    "F27" = 1346
    "F28" = 1356
    "F29" = 454
    "F32" = 162
    "F34" = 350
    "F35" = 59
    "F36" = 1017
    "F38" = 511
    "F39" = 666
    "F40" = 45
    "F41" = 45
    "F42" = 80
    "F45" = 166
    "F46" = 6
}

foreach ($cellRef in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($cellRef).Value = $exhibitionUpdates[$cellRef]
}

foreach ($cellRef in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range($cellRef).Value = $allTypesUpdates[$cellRef]
}
